$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "placeholder" FX rows (SOCIETE GENERALE / NOMURA HOLDINGS ticket
# stubs, rows 8-9) are no longer produced now that trade/repo/fx all come
# out of a single main.py run - drop them and let every row below shift up.
$ws.Range("A8:A9").EntireRow.Delete() | Out-Null

# Restore the cursor position left in the sheet when it was saved.
$ws.Range("C13").Select() | Out-Null

# Columns C (deal type), D (short name) and H (amount) were widened so the
# longer FX descriptions are no longer truncated.
$ws.Columns.Item(3).ColumnWidth = 16.833333333333332
$ws.Columns.Item(4).ColumnWidth = 18.833333333333332
$ws.Columns.Item(8).ColumnWidth = 14.666666666666666
